# Regenerate the "K" (strikeouts) column of save_data using the corrected
# per-game values (previously this column held a different "Strike#" metric).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new K value (column G), taken from the regenerated
# save_data. Row 20 already held the correct value and is left untouched.
$kValues = @{
    2  = 2
    3  = 4
    4  = 2
    5  = 1
    6  = 1
    7  = 1
    8  = 0
    9  = 0
    10 = 1
    11 = 0
    12 = 1
    13 = 1
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 1
    19 = 4
    21 = 2
    22 = 2
    23 = 2
    24 = 1
    25 = 3
    26 = 1
    27 = 2
    28 = 4
    29 = 1
    30 = 1
    31 = 2
    32 = 1
    33 = 1
    34 = 1
    35 = 1
    36 = 1
    37 = 0
    38 = 1
    39 = 1
    40 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
